$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the task text between B5, B6, B7 and update the difficulty styling
$ws.Range("B5").Value = "Create custom assets if needed"
$ws.Range("B5").Style = "Bad"

$ws.Range("B6").Value = "Creating tasklist"
$ws.Range("B6").Style = "Bad"

$ws.Range("B7").Value = "Model first level design in Unity"
$ws.Range("B7").Style = "Neutral"

# Update active cell / selection
$ws.Range("G15").Select()
